# Updates the cryptos list snapshot (Coin/Link/Price/Volume(1h)) to the
# latest scrape for the rows whose values changed. A couple of rows were
# also re-ranked (PEPE / WrappedeETH / Binance-PegBSC-USD around rows 28-30,
# and EthereumClassic / Monero around rows 39-40), so Coin + Link get
# rewritten there too.
#
# Price values are stored as plain text in the workbook (even when they look
# like plain numbers, e.g. "1.00"), so purely-numeric-looking price strings are
# written with a leading apostrophe (Excel's own text-literal marker) to stop
# them from being auto-converted to numbers -- exactly what typing them into
# Excel by hand would require.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '62.311.38'
$ws.Cells.Item(2, 5).Value = '  -2.22%  '

$ws.Cells.Item(3, 4).Value = '2.424.19'
$ws.Cells.Item(3, 5).Value = '  -2.54%  '

$ws.Cells.Item(4, 4).Value = '''1.00'
$ws.Cells.Item(4, 5).Value = '  +0.06%  '

$ws.Cells.Item(5, 4).Value = '''576.42'
$ws.Cells.Item(5, 5).Value = '  +0.07%  '

$ws.Cells.Item(6, 4).Value = '''142.46'
$ws.Cells.Item(6, 5).Value = '  -4.60%  '

$ws.Cells.Item(7, 5).Value = '  +0.19%  '

$ws.Cells.Item(8, 4).Value = '''0.525'
$ws.Cells.Item(8, 5).Value = '  -2.73%  '

$ws.Cells.Item(9, 4).Value = '2.419.96'
$ws.Cells.Item(9, 5).Value = '  -2.62%  '

$ws.Cells.Item(10, 4).Value = '''0.106'
$ws.Cells.Item(10, 5).Value = '  -5.68%  '

$ws.Cells.Item(11, 5).Value = '  +0.72%  '

$ws.Cells.Item(12, 4).Value = '''5.18'
$ws.Cells.Item(12, 5).Value = '  -1.63%  '

$ws.Cells.Item(13, 4).Value = '''0.345'
$ws.Cells.Item(13, 5).Value = '  -3.88%  '

$ws.Cells.Item(14, 4).Value = '''26.31'
$ws.Cells.Item(14, 5).Value = '  -3.26%  '

$ws.Cells.Item(15, 4).Value = '''0.0000172'
$ws.Cells.Item(15, 5).Value = '  -5.26%  '

$ws.Cells.Item(16, 4).Value = '2.881.79'
$ws.Cells.Item(16, 5).Value = '  -2.13%  '

$ws.Cells.Item(17, 4).Value = '62.341.49'
$ws.Cells.Item(17, 5).Value = '  -2.30%  '

$ws.Cells.Item(18, 4).Value = '2.424.93'
$ws.Cells.Item(18, 5).Value = '  -2.55%  '

$ws.Cells.Item(19, 4).Value = '''10.95'
$ws.Cells.Item(19, 5).Value = '  -5.33%  '

$ws.Cells.Item(20, 4).Value = '''7.07'
$ws.Cells.Item(20, 5).Value = '  -4.45%  '

$ws.Cells.Item(21, 4).Value = '''330.00'
$ws.Cells.Item(21, 5).Value = '  +0.03%  '

$ws.Cells.Item(22, 4).Value = '''4.11'
$ws.Cells.Item(22, 5).Value = '  -2.49%  '

$ws.Cells.Item(23, 4).Value = '''1.97'
$ws.Cells.Item(23, 5).Value = '  -6.08%  '

$ws.Cells.Item(24, 5).Value = '  +0.35%  '

$ws.Cells.Item(25, 4).Value = '''65.79'
$ws.Cells.Item(25, 5).Value = '  -0.67%  '

$ws.Cells.Item(26, 4).Value = '''628.41'
$ws.Cells.Item(26, 5).Value = '  -0.13%  '

$ws.Cells.Item(27, 4).Value = '''8.93'
$ws.Cells.Item(27, 5).Value = '  +3.42%  '

$ws.Cells.Item(28, 2).Value = 'WrappedeETH'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(28, 4).Value = '2.562.67'
$ws.Cells.Item(28, 5).Value = '  -4.82%  '

$ws.Cells.Item(29, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(29, 4).Value = '''1.00'
$ws.Cells.Item(29, 5).Value = '  +0.05%  '

$ws.Cells.Item(30, 2).Value = 'PEPE'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(30, 4).Value = '0.0₃0947'
$ws.Cells.Item(30, 5).Value = '  -9.55%  '

$ws.Cells.Item(31, 4).Value = '''1.42'
$ws.Cells.Item(31, 5).Value = '  -7.64%  '

$ws.Cells.Item(32, 4).Value = '''7.98'
$ws.Cells.Item(32, 5).Value = '  -5.15%  '

$ws.Cells.Item(33, 4).Value = '''1.88'
$ws.Cells.Item(33, 5).Value = '  -1.97%  '

$ws.Cells.Item(34, 5).Value = '  -3.14%  '

$ws.Cells.Item(35, 4).Value = '''4.95'
$ws.Cells.Item(35, 5).Value = '  -5.53%  '

$ws.Cells.Item(36, 5).Value = '  +0.46%  '

$ws.Cells.Item(37, 4).Value = '''1.43'
$ws.Cells.Item(37, 5).Value = '  -6.30%  '

$ws.Cells.Item(38, 4).Value = '''0.374'
$ws.Cells.Item(38, 5).Value = '  -2.89%  '

$ws.Cells.Item(39, 2).Value = 'Monero'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(39, 4).Value = '''149.02'
$ws.Cells.Item(39, 5).Value = '  +0.41%  '

$ws.Cells.Item(40, 2).Value = 'EthereumClassic'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(40, 4).Value = '''18.35'
$ws.Cells.Item(40, 5).Value = '  -2.58%  '

$ws.Cells.Item(41, 4).Value = '''5.24'
$ws.Cells.Item(41, 5).Value = '  -4.53%  '

$ws.Cells.Item(42, 4).Value = '''1.73'
$ws.Cells.Item(42, 5).Value = '  -5.28%  '

$ws.Cells.Item(43, 4).Value = '''42.42'
$ws.Cells.Item(43, 5).Value = '  +1.23%  '

$ws.Cells.Item(44, 5).Value = '  -0.01%  '

$ws.Cells.Item(45, 4).Value = '''2.45'
$ws.Cells.Item(45, 5).Value = '  -9.51%  '

$ws.Cells.Item(46, 4).Value = '''143.07'
$ws.Cells.Item(46, 5).Value = '  -4.80%  '

$ws.Cells.Item(47, 4).Value = '''3.65'
$ws.Cells.Item(47, 5).Value = '  -3.20%  '

$ws.Cells.Item(48, 4).Value = '''0.0517'
$ws.Cells.Item(48, 5).Value = '  -5.29%  '

$ws.Cells.Item(49, 4).Value = '''0.593'
$ws.Cells.Item(49, 5).Value = '  -2.29%  '

$ws.Cells.Item(50, 4).Value = '''19.32'
$ws.Cells.Item(50, 5).Value = '  -9.09%  '

$ws.Cells.Item(51, 4).Value = '0.0₆0237'
$ws.Cells.Item(51, 5).Value = '  +8.24%  '
